# Bercot ANF parser logic - time sheet update
# Rebuilds the "Date"/"Time" log: adds a new worked-time row (row 5, B5),
# plus two new rows (6 and 7) of actual clock in/out timestamps for 07/04,
# and widens the date/time columns to fit the new data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths (best effort; Excel quantizes ColumnWidth to whole pixels) ---
$ws.Range("A:B").ColumnWidth = 7.26
$ws.Range("C:C").ColumnWidth = 10.26
$ws.Range("D:D").ColumnWidth = 10.6

# --- Header row (values unchanged, format untouched) ---
$ws.Range("A1").Value = "Date"
$ws.Range("B1").Value = "Time"

# --- Existing rows 2-4 (dates/hours unchanged) ---
$ws.Range("A2").Value = 42915
$ws.Range("A2").NumberFormat = "MM/DD/YY"
$ws.Range("B2").Value = 1.5
$ws.Range("B2").NumberFormat = "General"

$ws.Range("A3").Value = 42917
$ws.Range("A3").NumberFormat = "MM/DD/YY"
$ws.Range("B3").Value = 3
$ws.Range("B3").NumberFormat = "General"
$ws.Range("C3").NumberFormat = "HH:MM:SS\ AM/PM"

$ws.Range("A4").Value = 42918
$ws.Range("A4").NumberFormat = "MM/DD/YY"
$ws.Range("B4").Value = 0.5
$ws.Range("B4").NumberFormat = "General"
$ws.Range("C4").NumberFormat = "HH:MM:SS\ AM/PM"

# --- Row 5: the old lone C5 time entry becomes a B5 worked-hours entry;
#     C5/D5 become blank (but time-formatted) clock in/out placeholders ---
$ws.Range("A5").Value = 42918
$ws.Range("A5").NumberFormat = "MM/DD/YY"
$ws.Range("B5").Value = 2.5
$ws.Range("B5").NumberFormat = "General"
$ws.Range("C5").ClearContents()
$ws.Range("C5").NumberFormat = "HH:MM:SS\ AM/PM"
$ws.Range("D5").NumberFormat = "HH:MM:SS\ AM/PM"

# --- Row 6 (new): 07/04 clock in/out pair ---
$ws.Range("A6").Value = 42920
$ws.Range("A6").NumberFormat = "MM/DD/YY"
$ws.Range("C6").Value = 0.28125
$ws.Range("C6").NumberFormat = "HH:MM:SS\ AM/PM"
$ws.Range("D6").Value = 0.338888888888889
$ws.Range("D6").NumberFormat = "HH:MM:SS\ AM/PM"

# --- Row 7 (new): 07/04 second clock-in ---
$ws.Range("A7").Value = 42920
$ws.Range("A7").NumberFormat = "MM/DD/YY"
$ws.Range("C7").Value = 0.368055555555556
$ws.Range("C7").NumberFormat = "HH:MM:SS\ AM/PM"

# --- Move active selection to A8, matching the post-edit cursor position ---
$ws.Range("A8").Select()
